# Fixed a bug in MergeSymbol
# The data rows (2-23) across columns A:F get reshuffled - each destination
# row receives the full A:F content that used to live at a different source
# row. Capture the original values first, then write them back out in the
# new order so that no data is lost while overwriting.

$ws = $excel.ActiveWorkbook.ActiveSheet

# Mapping: destination row -> source row (both reference the *original* sheet)
$rowMap = @{
    2  = 9
    3  = 12
    4  = 14
    5  = 4
    6  = 2
    7  = 7
    8  = 6
    9  = 11
    10 = 10
    11 = 5
    12 = 8
    13 = 3
    14 = 13
    15 = 15
    16 = 21
    17 = 16
    18 = 20
    19 = 17
    20 = 18
    21 = 19
    22 = 23
    23 = 22
}

# Snapshot original values for rows 2-23, columns A-F (1-6)
$original = @{}
foreach ($r in 2..23) {
    $rowVals = @{}
    foreach ($c in 1..6) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $original[$r] = $rowVals
}

# Write back values according to the mapping
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $original[$srcRow]
    foreach ($c in 1..6) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c]
    }
}
